$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook-level: absolute path + selection handled separately (absPath is environment-set, not scriptable) ---

# --- String cells (new shared strings get created in this left-to-right, top-to-bottom order) ---
$ws.Range("B35").Value = "ponta contorno"
$ws.Range("G35").Value = "detalhe 1"
$ws.Range("L35").Value = "detalhe 3"
$ws.Range("Q35").Value = "detalhe 5"
$ws.Range("G39").Value = "detalhe 2"
$ws.Range("L39").Value = "detalhe 4"
$ws.Range("B45").Value = "meio contorno"
$ws.Range("B56").Value = "abaixo meio contorno"
$ws.Range("B63").Value = "final contorno"

# --- Literal numeric cells ---
$ws.Range("O4").Value = 2352.9411764705883
$ws.Range("P4").Value = -9705.882352941177
$ws.Range("O5").Value = 705.8823529411765
$ws.Range("P5").Value = -9411.764705882353
$ws.Range("O6").Value = 1117.6470588235293
$ws.Range("P6").Value = -5000.0
$ws.Range("B36").Value = 0.0
$ws.Range("C36").Value = 0.9
$ws.Range("G36").Value = 0.07647058823529412
$ws.Range("H36").Value = 0.7058823529411765
$ws.Range("L36").Value = 0.11764705882352941
$ws.Range("M36").Value = -0.39999999999999997
$ws.Range("Q36").Value = 0.07058823529411765
$ws.Range("R36").Value = -0.9411764705882353
$ws.Range("B37").Value = 0.058823529411764705
$ws.Range("C37").Value = 0.8235294117647058
$ws.Range("G37").Value = 0.0
$ws.Range("H37").Value = 0.7058823529411765
$ws.Range("L37").Value = 0.0
$ws.Range("M37").Value = -0.39999999999999997
$ws.Range("Q37").Value = 0.0
$ws.Range("R37").Value = -0.9411764705882353
$ws.Range("B38").Value = 0.07647058823529412
$ws.Range("C38").Value = 0.7058823529411765
$ws.Range("D38").Value = 764.7058823529412
$ws.Range("E38").Value = 7058.823529411765
$ws.Range("B39").Value = 0.08823529411764706
$ws.Range("C39").Value = 0.5882352941176471
$ws.Range("D39").Value = 882.3529411764706
$ws.Range("E39").Value = 5882.352941176471
$ws.Range("B40").Value = -0.0882352941176471
$ws.Range("C40").Value = 0.5882352941176471
$ws.Range("D40").Value = -882.3529411764711
$ws.Range("E40").Value = 5882.352941176471
$ws.Range("G40").Value = 0.10588235294117647
$ws.Range("H40").Value = 0.4411764705882353
$ws.Range("L40").Value = 0.10588235294117647
$ws.Range("M40").Value = -0.49411764705882355
$ws.Range("B41").Value = -0.0764705882352941
$ws.Range("C41").Value = 0.7058823529411765
$ws.Range("D41").Value = -764.705882352941
$ws.Range("E41").Value = 7058.823529411765
$ws.Range("G41").Value = 0.0
$ws.Range("H41").Value = 0.4411764705882353
$ws.Range("L41").Value = 0.0
$ws.Range("M41").Value = -0.49411764705882355
$ws.Range("B42").Value = -0.0588235294117647
$ws.Range("C42").Value = 0.8235294117647058
$ws.Range("D42").Value = -588.235294117647
$ws.Range("E42").Value = 8235.294117647058
$ws.Range("B43").Value = 0.0
$ws.Range("C43").Value = 0.9
$ws.Range("D43").Value = 0.0
$ws.Range("E43").Value = 9000.0
$ws.Range("B46").Value = 0.08823529411764706
$ws.Range("C46").Value = 0.5882352941176471
$ws.Range("B47").Value = 0.10588235294117647
$ws.Range("C47").Value = 0.4411764705882353
$ws.Range("B48").Value = 0.11764705882352941
$ws.Range("C48").Value = -0.39999999999999997
$ws.Range("D48").Value = 1176.4705882352941
$ws.Range("E48").Value = -3999.9999999999995
$ws.Range("B49").Value = 0.10588235294117647
$ws.Range("C49").Value = -0.49411764705882355
$ws.Range("D49").Value = 1058.8235294117646
$ws.Range("E49").Value = -4941.176470588235
$ws.Range("B50").Value = -0.105882352941176
$ws.Range("C50").Value = -0.49411764705882355
$ws.Range("D50").Value = -1058.8235294117599
$ws.Range("E50").Value = -4941.176470588235
$ws.Range("B51").Value = -0.117647058823529
$ws.Range("C51").Value = -0.39999999999999997
$ws.Range("D51").Value = -1176.47058823529
$ws.Range("E51").Value = -3999.9999999999995
$ws.Range("B52").Value = -0.105882352941176
$ws.Range("C52").Value = 0.4411764705882353
$ws.Range("D52").Value = -1058.8235294117599
$ws.Range("E52").Value = 4411.764705882353
$ws.Range("B53").Value = -0.0882352941176471
$ws.Range("C53").Value = 0.5882352941176471
$ws.Range("D53").Value = -882.3529411764711
$ws.Range("E53").Value = 5882.352941176471
$ws.Range("B54").Value = 0.08823529411764706
$ws.Range("C54").Value = 0.5882352941176471
$ws.Range("D54").Value = 882.3529411764706
$ws.Range("E54").Value = 5882.352941176471
$ws.Range("B57").Value = 0.10588235294117647
$ws.Range("C57").Value = -0.49411764705882355
$ws.Range("B58").Value = 0.07058823529411765
$ws.Range("C58").Value = -0.9411764705882353
$ws.Range("B59").Value = -0.0705882352941176
$ws.Range("C59").Value = -0.9411764705882353
$ws.Range("D59").Value = -705.882352941176
$ws.Range("E59").Value = -9411.764705882353
$ws.Range("B60").Value = -0.105882352941176
$ws.Range("C60").Value = -0.49411764705882355
$ws.Range("D60").Value = -1058.8235294117599
$ws.Range("E60").Value = -4941.176470588235
$ws.Range("B61").Value = 0.10588235294117647
$ws.Range("C61").Value = -0.49411764705882355
$ws.Range("D61").Value = 1058.8235294117646
$ws.Range("E61").Value = -4941.176470588235
$ws.Range("B64").Value = 0.07058823529411765
$ws.Range("C64").Value = -0.9411764705882353
$ws.Range("B65").Value = 0.07058823529411765
$ws.Range("C65").Value = -1.0
$ws.Range("B66").Value = -0.0705882352941176
$ws.Range("C66").Value = -1.0
$ws.Range("D66").Value = -705.882352941176
$ws.Range("E66").Value = -10000.0
$ws.Range("B67").Value = -0.0705882352941176
$ws.Range("C67").Value = -0.9411764705882353
$ws.Range("D67").Value = -705.882352941176
$ws.Range("E67").Value = -9411.764705882353
$ws.Range("B68").Value = 0.07058823529411765
$ws.Range("C68").Value = -0.9411764705882353
$ws.Range("D68").Value = 705.8823529411765
$ws.Range("E68").Value = -9411.764705882353

# --- Formula cells / ranges (contiguous vertical shared-formula groups where applicable) ---
$ws.Range("O2").Formula = "=10000*M2"
$ws.Range("P2").Formula = "=10000*N2"
$ws.Range("O3:O5").Formula = "=10000*M3"
$ws.Range("P3:P5").Formula = "=10000*N3"
$ws.Range("O18").Formula = "=10000*M18"
$ws.Range("P18").Formula = "=10000*N18"
$ws.Range("O19").Formula = "=10000*M19"
$ws.Range("P19").Formula = "=10000*N19"
$ws.Range("D36").Formula = "=10000*B36"
$ws.Range("E36").Formula = "=10000*C36"
$ws.Range("I36").Formula = "=10000*G36"
$ws.Range("J36").Formula = "=10000*H36"
$ws.Range("N36").Formula = "=10000*L36"
$ws.Range("O36").Formula = "=10000*M36"
$ws.Range("S36").Formula = "=10000*Q36"
$ws.Range("T36").Formula = "=10000*R36"
$ws.Range("D37:D43").Formula = "=10000*B37"
$ws.Range("E37:E43").Formula = "=10000*C37"
$ws.Range("I37").Formula = "=10000*G37"
$ws.Range("J37").Formula = "=10000*H37"
$ws.Range("N37").Formula = "=10000*L37"
$ws.Range("O37").Formula = "=10000*M37"
$ws.Range("S37").Formula = "=10000*Q37"
$ws.Range("T37").Formula = "=10000*R37"
$ws.Range("I40").Formula = "=10000*G40"
$ws.Range("J40").Formula = "=10000*H40"
$ws.Range("N40").Formula = "=10000*L40"
$ws.Range("O40").Formula = "=10000*M40"
$ws.Range("I41").Formula = "=10000*G41"
$ws.Range("J41").Formula = "=10000*H41"
$ws.Range("N41").Formula = "=10000*L41"
$ws.Range("O41").Formula = "=10000*M41"
$ws.Range("D46").Formula = "=10000*B46"
$ws.Range("E46").Formula = "=10000*C46"
$ws.Range("D47:D54").Formula = "=10000*B47"
$ws.Range("E47:E54").Formula = "=10000*C47"
$ws.Range("D57").Formula = "=10000*B57"
$ws.Range("E57").Formula = "=10000*C57"
$ws.Range("D58:D61").Formula = "=10000*B58"
$ws.Range("E58:E61").Formula = "=10000*C58"
$ws.Range("D64").Formula = "=10000*B64"
$ws.Range("E64").Formula = "=10000*C64"
$ws.Range("D65:D68").Formula = "=10000*B65"
$ws.Range("E65:E68").Formula = "=10000*C65"

# --- Selection to match final workbook state ---
$ws.Range("S36:T37").Select() | Out-Null

